# Edit script reproducing the authored commit:
#   1. The table on slide 16 gets switched to a different built-in table
#      style (GUID {0AA76CFC-1CCC-446E-9D9E-05C023808280}) instead of the
#      previous custom style {3EEA1D48-8B3E-403A-9461-450E68F896A3}.
#   2. The deck's theme colour scheme (the "Integral" palette used by
#      ppt/theme/theme1.xml, the slide master's theme) is swapped for the
#      stock default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{0AA76CFC-1CCC-446E-9D9E-05C023808280}")
    }
}

# --- 2. Swap the theme's colour scheme from Integral -> Office -------
# RGB() packs as r + g*256 + b*65536 (PowerPoint/VBA colour order), so the
# integers below correspond to the stock Office theme hex values:
#   dk1 000000, lt1 FFFFFF, dk2 44546A, lt2 E7E6E6,
#   accent1 5B9BD5, accent2 ED7D31, accent3 A5A5A5, accent4 FFC000,
#   accent5 4472C4, accent6 70AD47, hlink 0563C1, folHlink 954F72
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Range().ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
